$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

Replace-Text "2025-12-16 Tuesday" "2025-12-17 Wednesday"

Replace-Text "521×3=" "872×8="
Replace-Text "124×4=" "649×4="
Replace-Text "638×5=" "809×7="
Replace-Text "651×4=" "182×4="
Replace-Text "656×2=" "984×2="

Replace-Text "970×5=" "246×4="
Replace-Text "827×9=" "329×7="
Replace-Text "948×9=" "626×9="
Replace-Text "540×7=" "868×4="
Replace-Text "579×3=" "947×7="

Replace-Text "115×7=" "392×8="
Replace-Text "683×2=" "324×8="
Replace-Text "595×9=" "577×5="
Replace-Text "618×7=" "915×4="
Replace-Text "560×2=" "743×7="

Replace-Text "375×9=" "873×8="
Replace-Text "691×5=" "605×9="
Replace-Text "858×5=" "716×2="
Replace-Text "312×4=" "321×9="
Replace-Text "807×7=" "815×9="

Replace-Text "889×9=" "164×4="
Replace-Text "666×8=" "381×8="
Replace-Text "473×8=" "690×5="
Replace-Text "716×8=" "973×7="
Replace-Text "649×3=" "657×8="
